# Update the Schedule (Hoja1) planning values for team/individual plan tracking.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Week 1 (row 14): direct/cumulative hours revised down from 4 to 3.5
$ws.Range("D14").Value = 3.5
$ws.Range("E14").Value = 3.5

# Week 3 (row 16): direct hours revised from 6 to 9.2
$ws.Range("D16").Value = 9.2

# Week 4 (row 17): direct hours revised from 6 to 4
$ws.Range("D17").Value = 4

# Week 5 (row 18): direct hours revised from 7 to 11
$ws.Range("D18").Value = 11

# Week 6 (row 19): direct hours revised from 7 to 2.5
$ws.Range("D19").Value = 2.5

# Move the active selection to H22 (previously G16:H16)
$ws.Activate() | Out-Null
$ws.Range("H22").Select() | Out-Null
